# LWA-352 Antenna Status & System Configuration
# "pretest change for demo"
#
# Semantic edit: antenna LWA-213's SNAP2 chassis # (K4) changes from 3 -> 2.
#
# The source XML diff also shows every column on Sheet1 getting a couple of
# pixels narrower (a side effect of the file having been re-saved by a
# different Excel/Calc build). We reproduce that as closely as the
# ColumnWidth object model here allows, column-range by column-range, to
# keep the sheet's visual layout in step with the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (narrowed slightly across the board) ---------------
$ws.Range("A1").EntireColumn.ColumnWidth = 8.212585034013607
$ws.Range("B1").EntireColumn.ColumnWidth = 9.426870748299367
$ws.Range("C1").EntireColumn.ColumnWidth = 12.937074829931966
$ws.Range("D1").EntireColumn.ColumnWidth = 14.554421768707465
$ws.Range("E1").EntireColumn.ColumnWidth = 12.666666666666666
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 9.289115646258466
$ws.Range("H1").EntireColumn.ColumnWidth = 11.452380952380967
$ws.Range("I1").EntireColumn.ColumnWidth = 21.845238095238066
$ws.Range("J1").EntireColumn.ColumnWidth = 8.345238095238097
$ws.Range("K1").EntireColumn.ColumnWidth = 9.697278911564666
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 10.508503401360565
$ws.Range("N1").EntireColumn.ColumnWidth = 11.722789115646266
$ws.Range("O1").EntireColumn.ColumnWidth = 7.130952380952377
$ws.Range("P1").EntireColumn.ColumnWidth = 17.661564625850367
$ws.Range("Q1").EntireColumn.ColumnWidth = 11.855442176870767
$ws.Range("R1:S1").EntireColumn.ColumnWidth = 11.044217687074866
$ws.Range("T1").EntireColumn.ColumnWidth = 9.156462585034015
$ws.Range("U1:V1").EntireColumn.ColumnWidth = 9.559523809523766
$ws.Range("W1").EntireColumn.ColumnWidth = 11.722789115646266
$ws.Range("X1").EntireColumn.ColumnWidth = 11.452380952380967
$ws.Range("Y1:Z1").EntireColumn.ColumnWidth = 11.722789115646266
$ws.Range("AA1").EntireColumn.ColumnWidth = 35.74829931972786

# --- Cell edit: LWA-213 SNAP2 chassis # 3 -> 2 -------------------------
$ws.Range("K4").Value = 2
